$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.201.91'
$ws.Range("E2").Value = '  -7.15%  '
$ws.Range("D3").Value = '1.673.23'
$ws.Range("E3").Value = '  -4.40%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''218.50'
$ws.Range("E5").Value = '  -4.21%  '
$ws.Range("D6").Value = '''0.5103'
$ws.Range("E6").Value = '  -12.55%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").Value = '''21.84'
$ws.Range("E9").Value = '  -4.98%  '
$ws.Range("E10").Value = '  -4.58%  '
$ws.Range("D11").Value = '''0.07385'
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("D12").Value = '1.672.83'
$ws.Range("E12").Value = '  -4.70%  '
$ws.Range("D13").Value = '''4.544'
$ws.Range("E13").Value = '  -3.69%  '
$ws.Range("D14").Value = '''0.5754'
$ws.Range("E14").Value = '  -4.62%  '
$ws.Range("D15").Value = '1.903.83'
$ws.Range("E15").Value = '  -4.20%  '
$ws.Range("D16").Value = '''0.000008540'
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '''64.99'
$ws.Range("E17").Value = '  -12.33%  '
$ws.Range("D18").Value = '26.295.90'
$ws.Range("E18").Value = '  -6.74%  '
$ws.Range("E19").Value = '  -5.73%  '
$ws.Range("D20").Value = '''1.006'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").Value = '''10.86'
$ws.Range("E21").Value = '  -3.78%  '
$ws.Range("D22").Value = '''187.08'
$ws.Range("E22").Value = '  -8.92%  '
$ws.Range("D23").Value = '''6.215'
$ws.Range("E23").Value = '  -6.89%  '
$ws.Range("D24").Value = '''1.007'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '''143.53'
$ws.Range("E25").Value = '  -4.87%  '
$ws.Range("D26").Value = '''7.565'
$ws.Range("E26").Value = '  -5.80%  '
$ws.Range("D27").Value = '''0.1173'
$ws.Range("E27").Value = '  -5.09%  '
$ws.Range("D28").Value = '''15.70'
$ws.Range("E28").Value = '  -2.38%  '
$ws.Range("D29").Value = '''1.324'
$ws.Range("E29").Value = '  -6.06%  '
$ws.Range("D30").Value = '''0.05795'
$ws.Range("E30").Value = '  -5.26%  '
$ws.Range("D31").Value = '''1.331'
$ws.Range("E31").Value = '  -4.92%  '
$ws.Range("D32").Value = '''3.507'
$ws.Range("E32").Value = '  -5.96%  '
$ws.Range("D33").Value = '''3.501'
$ws.Range("E33").Value = '  -5.60%  '
$ws.Range("D34").Value = '''1.665'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = '''1.002'
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("D36").Value = '''0.5981'
$ws.Range("E36").Value = '  -5.42%  '
$ws.Range("D37").Value = '''2.370'
$ws.Range("E37").Value = '  +1.53%  '
$ws.Range("D38").Value = '''2.647'
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").Value = '1.098.13'
$ws.Range("E39").Value = '  -3.04%  '
$ws.Range("D40").Value = '''0.01602'
$ws.Range("E40").Value = '  -3.66%  '
$ws.Range("D41").Value = '''5.913'
$ws.Range("E41").Value = '  -5.89%  '
$ws.Range("D42").Value = '''0.8629'
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").Value = '''1.006'
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").Value = '''99.44'
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = '1.821.80'
$ws.Range("E45").Value = '  -4.13%  '
$ws.Range("D46").Value = '''0.00000000114'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("D47").Value = '''56.25'
$ws.Range("E47").Value = '  -4.70%  '
$ws.Range("D48").Value = '''1.006'
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("D49").Value = '''8.045'
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("E50").Value = '  -3.33%  '
$ws.Range("D51").Value = '''0.05207'
